$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("matmul-optimization")
$ws.Range("K3").Value = "Test"
